{"js": "// Update each division-problem cell in the table to its new value.\nconst replacements = [\n  [\"509\u00f77=\", \"428\u00f73=\"],\n  [\"101\u00f78=\", \"956\u00f75=\"],\n  [\"756\u00f78=\", \"314\u00f77=\"],\n  [\"504\u00f74=\", \"922\u00f75=\"],\n  [\"111\u00f79=\", \"271\u00f79=\"],\n  [\"601\u00f77=\", \"787\u00f78=\"],\n  [\"389\u00f72=\", \"564\u00f76=\"],\n  [\"997\u00f73=\", \"153\u00f79=\"],\n  [\"335\u00f78=\", \"576\u00f73=\"],\n  [\"198\u00f75=\", \"903\u00f76=\"],\n  [\"985\u00f78=\", \"480\u00f77=\"],\n  [\"185\u00f76=\", \"630\u00f77=\"],\n  [\"751\u00f76=\", \"733\u00f75=\"],\n  [\"350\u00f76=\", \"118\u00f78=\"],\n  [\"455\u00f76=\", \"923\u00f77=\"],\n  [\"609\u00f78=\", \"213\u00f77=\"],\n  [\"365\u00f73=\", \"366\u00f76=\"],\n  [\"427\u00f76=\", \"416\u00f79=\"],\n  [\"871\u00f73=\", \"858\u00f76=\"],\n  [\"173\u00f77=\", \"750\u00f77=\"],\n  [\"292\u00f73=\", \"669\u00f79=\"],\n  [\"502\u00f73=\", \"501\u00f77=\"],\n  [\"102\u00f74=\", \"856\u00f79=\"],\n  [\"732\u00f76=\", \"109\u00f78=\"],\n  [\"392\u00f76=\", \"129\u00f72=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update each division-problem cell in the table to the new value.\n# wdFindContinue = 1, wdReplaceAll = 2\n$d = $word.ActiveDocument\n\n$null = $d.Content.Find.Execute(\"509\u00f77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"428\u00f73=\", 2)\n$null = $d.Content.Find.Execute(\"101\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"956\u00f75=\", 2)\n$null = $d.Content.Find.Execute(\"756\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"314\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"504\u00f74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"922\u00f75=\", 2)\n$null = $d.Content.Find.Execute(\"111\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"271\u00f79=\", 2)\n$null = $d.Content.Find.Execute(\"601\u00f77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"787\u00f78=\", 2)\n$null = $d.Content.Find.Execute(\"389\u00f72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"564\u00f76=\", 2)\n$null = $d.Content.Find.Execute(\"997\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"153\u00f79=\", 2)\n$null = $d.Content.Find.Execute(\"335\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"576\u00f73=\", 2)\n$null = $d.Content.Find.Execute(\"198\u00f75=\", $false, $false, $false, $false, $false, $true, 1, $false, \"903\u00f76=\", 2)\n$null = $d.Content.Find.Execute(\"985\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"480\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"185\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"630\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"751\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"733\u00f75=\", 2)\n$null = $d.Content.Find.Execute(\"350\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"118\u00f78=\", 2)\n$null = $d.Content.Find.Execute(\"455\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"923\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"609\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"213\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"365\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"366\u00f76=\", 2)\n$null = $d.Content.Find.Execute(\"427\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"416\u00f79=\", 2)\n$null = $d.Content.Find.Execute(\"871\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"858\u00f76=\", 2)\n$null = $d.Content.Find.Execute(\"173\u00f77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"750\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"292\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"669\u00f79=\", 2)\n$null = $d.Content.Find.Execute(\"502\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"501\u00f77=\", 2)\n$null = $d.Content.Find.Execute(\"102\u00f74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"856\u00f79=\", 2)\n$null = $d.Content.Find.Execute(\"732\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"109\u00f78=\", 2)\n$null = $d.Content.Find.Execute(\"392\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"129\u00f72=\", 2)\n"}
